# Reorder the AF:AJ "Usage" columns on row 1 (headers) and row 2 (values).
# The underlying label -> value mapping stays the same:
#   pkmUsage      = 858
#   tkm-N3Usage   = 130.3
#   tkm-SZMUsage  = 414.5
#   tkm-N1Usage   = 7.5
#   tkm-N2Usage   = 24.2
# but the columns they live in (and the shared-string table order) change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header order for columns AF1..AJ1
$ws.Range("AF1").Value = "tkm-N3Usage"
$ws.Range("AG1").Value = "tkm-N2Usage"
$ws.Range("AH1").Value = "tkm-SZMUsage"
$ws.Range("AI1").Value = "pkmUsage"
$ws.Range("AJ1").Value = "tkm-N1Usage"

# New values matching the headers above
$ws.Range("AF2").Value = 130.3
$ws.Range("AG2").Value = 24.2
$ws.Range("AH2").Value = 414.5
$ws.Range("AI2").Value = 858
$ws.Range("AJ2").Value = 7.5
